$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 2: replace the old article (790-EADOGM204SA) with the new connector part
$ws.Range("A2").Value = "200-FLE10501GDVKTR"
$ws.Range("B2").Value = "Platine-zu-Platine & Mezzanine-Steckverbinder Cost Effective Surface Mount Socket, 0.050`" Pitch"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "Stck"
$ws.Range("E2").Value = 2.3
$ws.Range("H2").Value = "Mouser"

# Row 3: fill in the previously-blank (formatting-only) row with a new article
$ws.Range("A3").Value = "200-FTSH10501LDVK"
$ws.Range("B3").Value = "Sockel & Kabelgehäuse High Reliability Header Strips, .050`" pitch"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "Stck"
$ws.Range("E3").Value = 2.36
# H3 previously carried a stray explicit black-font style; pull the (unstyled)
# format from H2 before writing the value so it reverts to the default style.
$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = "Mouser"

# Row 4: another new article; give D4 the same "dropdown" style used by D2/D3
$ws.Range("A4").Value = "667-ERA-3AEB2800V"
$ws.Range("B4").Value = "Dünnfilmwiderstände - SMD 0603 280ohm 0.1% 25ppm"
$ws.Range("C4").Value = 20
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "Stck"
$ws.Range("E4").Value = 0.075
$ws.Range("H4").Value = "Mouser"

$excel.Application.CutCopyMode = $false

$ws.Range("H5").Select()
